$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# E2: add new value
$ws.Range("E2").Value = "26 TL - 26 TL"

# D6: clear existing value
$ws.Range("D6").Value = ""

# D12: clear existing value
$ws.Range("D12").Value = ""

# C13: update value
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"

# D13: clear existing value
$ws.Range("D13").Value = ""

# K13: update value
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# D14: clear existing value
$ws.Range("D14").Value = ""
